$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1016.7619
$ws.Range("J40").Value = 852
$ws.Range("L40").Value = 852
$ws.Range("N40").Value = -1202

$ws.Range("H62").Value = 7964.143
$ws.Range("I62").Value = 6437.25
$ws.Range("K62").Value = 6437.25
$ws.Range("M62").Value = -5813.25

$ws.Range("H65").Value = 7964.143
$ws.Range("I65").Value = 6437.25
$ws.Range("K65").Value = 32186.25
$ws.Range("M65").Value = -29066.25

$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 3380
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 10140
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -10680

$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 3380
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 10140
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -12012

$ws.Range("H100").Value = 2488.4666
$ws.Range("J100").Value = 2807.8333
$ws.Range("L100").Value = 2807.8333
$ws.Range("N100").Value = -3889.8333

$ws.Range("H113").Value = 35128.547
$ws.Range("I113").Value = 52778.617
$ws.Range("J113").Value = 4240.9165
$ws.Range("K113").Value = 52778.617
$ws.Range("L113").Value = 4240.9165
$ws.Range("M113").Value = -49524.617
$ws.Range("N113").Value = -10748.9165

$ws.Range("H132").Value = 10480314
$ws.Range("I132").Value = 11529939
$ws.Range("K132").Value = 34589817
$ws.Range("M132").Value = -34587287

$ws.Range("H137").Value = 1551
$ws.Range("I137").Value = 1193.25
$ws.Range("J137").Value = 2624.25
$ws.Range("K137").Value = 3579.75
$ws.Range("L137").Value = 7872.75
$ws.Range("M137").Value = -1029.75
$ws.Range("N137").Value = -12972.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2146.6316
$ws.Range("I2").Value = 2656.9092
$ws.Range("J2").Value = 1445
$ws.Range("K2").Value = 2656.9092
$ws.Range("L2").Value = 1445
$ws.Range("M2").Value = -2543.9092
$ws.Range("N2").Value = -1671

$ws.Range("H45").Value = 11678.519
$ws.Range("I45").Value = 9228.857
$ws.Range("K45").Value = 9228.857
$ws.Range("M45").Value = -8851.857

$ws.Range("H102").Value = 4959.0835
$ws.Range("I102").Value = 5418.1665
$ws.Range("K102").Value = 5418.1665
$ws.Range("M102").Value = -3796.1665

$ws.Range("H110").Value = 6421.1177
$ws.Range("I110").Value = 7087.2104
$ws.Range("J110").Value = 5577.4
$ws.Range("K110").Value = 7087.2104
$ws.Range("L110").Value = 5577.4
$ws.Range("M110").Value = -5042.2104
$ws.Range("N110").Value = -9667.4

$ws.Range("H116").Value = 2146.6316
$ws.Range("I116").Value = 2656.9092
$ws.Range("J116").Value = 1445
$ws.Range("K116").Value = 2656.9092
$ws.Range("L116").Value = 1445
$ws.Range("M116").Value = -362.9092000000001
$ws.Range("N116").Value = -6033

$ws.Range("H122").Value = 2704.875
$ws.Range("I122").Value = 2704.875
$ws.Range("K122").Value = 8114.625
$ws.Range("M122").Value = -5664.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2146.6316
$ws.Range("I3").Value = 2656.9092
$ws.Range("J3").Value = 1445
$ws.Range("K3").Value = 2656.9092
$ws.Range("L3").Value = 1445
$ws.Range("M3").Value = -2542.9092
$ws.Range("N3").Value = -1673

$ws.Range("H6").Value = 40427
$ws.Range("J6").Value = 40427
$ws.Range("L6").Value = 40427
$ws.Range("N6").Value = -40653

$ws.Range("H86").Value = 14368
$ws.Range("I86").Value = 13376.667
$ws.Range("J86").Value = 18333.334
$ws.Range("K86").Value = 13376.667
$ws.Range("L86").Value = 18333.334
$ws.Range("M86").Value = -12253.667
$ws.Range("N86").Value = -20579.334

$ws.Range("H89").Value = 14368
$ws.Range("I89").Value = 13376.667
$ws.Range("J89").Value = 18333.334
$ws.Range("K89").Value = 66883.33499999999
$ws.Range("L89").Value = 91666.67
$ws.Range("M89").Value = -61267.33499999999
$ws.Range("N89").Value = -102898.67

$ws.Range("H105").Value = 2799.8965
$ws.Range("I105").Value = 2813.15
$ws.Range("J105").Value = 2770.4443
$ws.Range("K105").Value = 2813.15
$ws.Range("L105").Value = 2770.4443
$ws.Range("M105").Value = -1066.15
$ws.Range("N105").Value = -6264.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 8674.200000000001
$ws.Range("I55").Value = 8674.200000000001
$ws.Range("K55").Value = 8674.200000000001
$ws.Range("M55").Value = -8359.200000000001

$ws.Range("H70").Value = 40000
$ws.Range("J70").Value = 40000
$ws.Range("L70").Value = 40000
$ws.Range("N70").Value = -40630

$ws.Range("H73").Value = 40000
$ws.Range("J73").Value = 40000
$ws.Range("L73").Value = 40000
$ws.Range("N73").Value = -42184

$ws.Range("H74").Value = 49750
$ws.Range("J74").Value = 49750
$ws.Range("L74").Value = 49750
$ws.Range("N74").Value = -51498

$ws.Range("H77").Value = 49750
$ws.Range("J77").Value = 49750
$ws.Range("L77").Value = 149250
$ws.Range("N77").Value = -157986

$ws.Range("H97").Value = 17870
$ws.Range("J97").Value = 17870
$ws.Range("L97").Value = 17870
$ws.Range("N97").Value = -19852

$ws.Range("H105").Value = 1509.4
$ws.Range("I105").Value = 1392.5333
$ws.Range("J105").Value = 1860
$ws.Range("K105").Value = 1392.5333
$ws.Range("L105").Value = 1860
$ws.Range("M105").Value = 354.4666999999999
$ws.Range("N105").Value = -5354

$ws.Range("H134").Value = 19138.963
$ws.Range("I134").Value = 11775.842
$ws.Range("K134").Value = 35327.526
$ws.Range("M134").Value = -32792.526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2079.75
$ws.Range("I81").Value = 512.5
$ws.Range("J81").Value = 3647
$ws.Range("K81").Value = 1537.5
$ws.Range("L81").Value = 10941
$ws.Range("M81").Value = -414.5
$ws.Range("N81").Value = -13187

$ws.Range("H84").Value = 2079.75
$ws.Range("I84").Value = 512.5
$ws.Range("J84").Value = 3647
$ws.Range("K84").Value = 4612.5
$ws.Range("L84").Value = 32823
$ws.Range("M84").Value = 1003.5
$ws.Range("N84").Value = -44055

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 234.70833
$ws.Range("J2").Value = 310.6
$ws.Range("L2").Value = 310.6
$ws.Range("N2").Value = -536.6

$ws.Range("H80").Value = 3497
$ws.Range("I80").Value = 3995
$ws.Range("J80").Value = 2999
$ws.Range("K80").Value = 3995
$ws.Range("L80").Value = 2999
$ws.Range("M80").Value = -2997
$ws.Range("N80").Value = -4995

$ws.Range("H83").Value = 3497
$ws.Range("I83").Value = 3995
$ws.Range("J83").Value = 2999
$ws.Range("K83").Value = 19975
$ws.Range("L83").Value = 14995
$ws.Range("M83").Value = -14983
$ws.Range("N83").Value = -24979

$ws.Range("H107").Value = 1770.3
$ws.Range("I107").Value = 1930
$ws.Range("K107").Value = 1930
$ws.Range("M107").Value = -10

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4593.5
$ws.Range("I16").Value = 1759.6923
$ws.Range("J16").Value = 7427.3076
$ws.Range("K16").Value = 1759.6923
$ws.Range("L16").Value = 7427.3076
$ws.Range("M16").Value = -1589.6923
$ws.Range("N16").Value = -7767.3076

$ws.Range("H22").Value = 2007
$ws.Range("J22").Value = 2007
$ws.Range("L22").Value = 2007
$ws.Range("N22").Value = -2597

$ws.Range("H27").Value = 2007
$ws.Range("J27").Value = 2007
$ws.Range("L27").Value = 2007
$ws.Range("N27").Value = -2221

$ws.Range("H48").Value = 5000
$ws.Range("I48").Value = 5000
$ws.Range("K48").Value = 5000
$ws.Range("M48").Value = -4339

$ws.Range("H61").Value = 2740.725
$ws.Range("I61").Value = 2509.3547
$ws.Range("J61").Value = 3537.6667
$ws.Range("K61").Value = 2509.3547
$ws.Range("L61").Value = 3537.6667
$ws.Range("M61").Value = -2307.3547
$ws.Range("N61").Value = -3941.6667

$ws.Range("H82").Value = 2740
$ws.Range("I82").Value = 1900
$ws.Range("J82").Value = 2950
$ws.Range("K82").Value = 1900
$ws.Range("L82").Value = 2950
$ws.Range("M82").Value = -1539
$ws.Range("N82").Value = -3672

$ws.Range("H85").Value = 2740
$ws.Range("I85").Value = 1900
$ws.Range("J85").Value = 2950
$ws.Range("K85").Value = 1900
$ws.Range("L85").Value = 2950
$ws.Range("M85").Value = -652
$ws.Range("N85").Value = -5446

$ws.Range("H99").Value = 19999
$ws.Range("I99").Value = 19999
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 19999
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -17004
$ws.Range("N99").ClearContents()

$ws.Range("H100").Value = 3442.5
$ws.Range("I100").Value = 3577.1428
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 3577.1428
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -3036.1428
$ws.Range("N100").Value = -3582

$ws.Range("H113").Value = 2740.725
$ws.Range("I113").Value = 2509.3547
$ws.Range("J113").Value = 3537.6667
$ws.Range("K113").Value = 2509.3547
$ws.Range("L113").Value = 3537.6667
$ws.Range("M113").Value = -339.3546999999999
$ws.Range("N113").Value = -7877.6667

$ws.Range("H132").Value = 3773.814
$ws.Range("I132").Value = 3391.52
$ws.Range("J132").Value = 4304.778
$ws.Range("K132").Value = 10174.56
$ws.Range("L132").Value = 12914.334
$ws.Range("M132").Value = -7644.559999999999
$ws.Range("N132").Value = -17974.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 37975
$ws.Range("J70").Value = 37975
$ws.Range("L70").Value = 37975
$ws.Range("N70").Value = -38605

$ws.Range("H73").Value = 37975
$ws.Range("J73").Value = 37975
$ws.Range("L73").Value = 37975
$ws.Range("N73").Value = -40159

$ws.Range("H75").Value = 36000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872

$ws.Range("H78").Value = 36000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360

$ws.Range("H122").Value = 3783.9092
$ws.Range("I122").Value = 3462.3
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 10386.9
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -7936.900000000001
$ws.Range("N122").Value = -25900
